$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.40597222990741
$ws.Range("C2").Value = 9.44288880431767
$ws.Range("D2").Value = 5.9820434192404
$ws.Range("E2").Value = 10.2293753987947
$ws.Range("G2").Value = 39.97330704477238
$ws.Range("H2").Value = 16.12124930084249
$ws.Range("M2").Value = 15.86734605725166
$ws.Range("N2").Value = 17.81794998644376
$ws.Range("B3").Value = 14.79867437197658
$ws.Range("C3").Value = 8.818381737864085
$ws.Range("D3").Value = 5.863929050597995
$ws.Range("E3").Value = 10.13407516315474
$ws.Range("G3").Value = 39.20635205070928
$ws.Range("H3").Value = 16.07967375491414
$ws.Range("M3").Value = 15.58392473444794
$ws.Range("N3").Value = 17.8867236944457
$ws.Range("B4").Value = 14.41742632606623
$ws.Range("C4").Value = 8.411130888197571
$ws.Range("D4").Value = 5.792176720368224
$ws.Range("E4").Value = 10.07846394674266
$ws.Range("G4").Value = 38.74399187591163
$ws.Range("H4").Value = 16.05857345168006
$ws.Range("M4").Value = 15.41226996672802
$ws.Range("N4").Value = 17.93091038690633
$ws.Range("B5").Value = 14.26022838013124
$ws.Range("C5").Value = 8.239169046025371
$ws.Range("D5").Value = 5.763174448529601
$ws.Range("E5").Value = 10.05655200598757
$ws.Range("G5").Value = 38.55800740966662
$ws.Range("H5").Value = 16.05109013152067
$ws.Range("M5").Value = 15.34301268598015
$ws.Range("N5").Value = 17.94941033091288
$ws.Range("B6").Value = 14.23402322144675
$ws.Range("C6").Value = 8.210251748831912
$ws.Range("D6").Value = 5.758374259061695
$ws.Range("E6").Value = 10.05295939705049
$ws.Range("G6").Value = 38.52727978097556
$ws.Range("H6").Value = 16.04991492095913
$ws.Range("M6").Value = 15.33155730723822
$ws.Range("N6").Value = 17.95251206336694
$ws.Range("B7").Value = 14.41531336139606
$ws.Range("C7").Value = 8.408836084735597
$ws.Range("D7").Value = 5.791784566602321
$ws.Range("E7").Value = 10.07816537251559
$ws.Range("G7").Value = 38.74147341352582
$ws.Range("H7").Value = 16.0584680120078
$ws.Range("M7").Value = 15.41133300183361
$ws.Range("N7").Value = 17.93115788442366
$ws.Range("B8").Value = 15.19846376828739
$ws.Range("C8").Value = 9.232487224976834
$ws.Range("D8").Value = 5.941183824480611
$ws.Range("E8").Value = 10.19592575792269
$ws.Range("G8").Value = 39.70725468841077
$ws.Range("H8").Value = 16.10599672537257
$ws.Range("M8").Value = 15.76919012110561
$ws.Range("N8").Value = 17.84125698567259
$ws.Range("B9").Value = 16.65749892949989
$ws.Range("C9").Value = 10.65967709453102
$ws.Range("D9").Value = 6.23836896801541
$ws.Range("E9").Value = 10.44893404069946
$ws.Range("G9").Value = 41.65663060966042
$ws.Range("H9").Value = 16.2341873068523
$ws.Range("M9").Value = 16.48541582593848
$ws.Range("N9").Value = 17.68046994769717
$ws.Range("B10").Value = 17.67067196961279
$ws.Range("C10").Value = 11.59483576307588
$ws.Range("D10").Value = 6.456797542847403
$ws.Range("E10").Value = 10.64689848728427
$ws.Range("G10").Value = 43.10622877007021
$ws.Range("H10").Value = 16.3494292409953
$ws.Range("M10").Value = 17.01474938148977
$ws.Range("N10").Value = 17.57173881121565
$ws.Range("B11").Value = 18.11678018811254
$ws.Range("C11").Value = 11.99589890191905
$ws.Range("D11").Value = 6.555686122629085
$ws.Range("E11").Value = 10.73926285945416
$ws.Range("G11").Value = 43.76610075498905
$ws.Range("H11").Value = 16.40634564623066
$ws.Range("M11").Value = 17.25507142917101
$ws.Range("N11").Value = 17.52430375522423
$ws.Range("B12").Value = 18.28344196227214
$ws.Range("C12").Value = 12.14429046385929
$ws.Range("D12").Value = 6.593025069982525
$ws.Range("E12").Value = 10.77454307159053
$ws.Range("G12").Value = 44.01577323665408
$ws.Range("H12").Value = 16.42853537950929
$ws.Range("M12").Value = 17.34591434000243
$ws.Range("N12").Value = 17.5066321274936
$ws.Range("B13").Value = 18.24765138393957
$ws.Range("C13").Value = 12.11248627047007
$ws.Range("D13").Value = 6.584988871912126
$ws.Range("E13").Value = 10.76693178878136
$ws.Range("G13").Value = 43.96201516892609
$ws.Range("H13").Value = 16.42372826452404
$ws.Range("M13").Value = 17.32635833796612
$ws.Range("N13").Value = 17.51042509438838
$ws.Range("B14").Value = 18.13053779096714
$ws.Range("C14").Value = 12.00817685075396
$ws.Range("D14").Value = 6.558760404480757
$ws.Range("E14").Value = 10.74215945812357
$ws.Range("G14").Value = 43.78664710850336
$ws.Range("H14").Value = 16.40815850515417
$ws.Range("M14").Value = 17.26254889506861
$ws.Range("N14").Value = 17.52284407060283
$ws.Range("B15").Value = 18.05850283870052
$ws.Range("C15").Value = 11.94383136271059
$ws.Range("D15").Value = 6.542679504452612
$ws.Range("E15").Value = 10.72702441077356
$ws.Range("G15").Value = 43.67919437021457
$ws.Range("H15").Value = 16.39870420391254
$ws.Range("M15").Value = 17.22343998501685
$ws.Range("N15").Value = 17.53048893380917
$ws.Range("B16").Value = 17.64120746761911
$ws.Range("C16").Value = 11.5681370352293
$ws.Range("D16").Value = 6.450321936889802
$ws.Range("E16").Value = 10.64090632540875
$ws.Range("G16").Value = 43.0630901063687
$ws.Range("H16").Value = 16.34579925137049
$ws.Range("M16").Value = 16.99902612969064
$ws.Range("N16").Value = 17.57487955233636
$ws.Range("B17").Value = 17.38131421958655
$ws.Range("C17").Value = 11.33143849064126
$ws.Range("D17").Value = 6.393512900272261
$ws.Range("E17").Value = 10.58864722520931
$ws.Range("G17").Value = 42.6850473593712
$ws.Range("H17").Value = 16.31448814125646
$ws.Range("M17").Value = 16.86116611988331
$ws.Range("N17").Value = 17.60263061191635
$ws.Range("B18").Value = 17.23044567519639
$ws.Range("C18").Value = 11.19300315690733
$ws.Range("D18").Value = 6.360795317400078
$ws.Range("E18").Value = 10.5588082729199
$ws.Range("G18").Value = 42.46766429703421
$ws.Range("H18").Value = 16.29690228297324
$ws.Range("M18").Value = 16.78183339353037
$ws.Range("N18").Value = 17.61878313402319
$ws.Range("B19").Value = 17.17913102038772
$ws.Range("C19").Value = 11.1457370135954
$ws.Range("D19").Value = 6.349711642301663
$ws.Range("E19").Value = 10.54874379891086
$ws.Range("G19").Value = 42.39408036985231
$ws.Range("H19").Value = 16.29102100668037
$ws.Range("M19").Value = 16.7549688642629
$ws.Range("N19").Value = 17.6242848950042
$ws.Range("B20").Value = 17.40912472245721
$ws.Range("C20").Value = 11.35687259529793
$ws.Range("D20").Value = 6.399565014127455
$ws.Range("E20").Value = 10.59418782263334
$ws.Range("G20").Value = 42.72528673920153
$ws.Range("H20").Value = 16.31777749882094
$ws.Range("M20").Value = 16.87584631110045
$ws.Range("N20").Value = 17.59965671625428
$ws.Range("B21").Value = 18.16499953466503
$ws.Range("C21").Value = 12.03890941866009
$ws.Range("D21").Value = 6.566467577895456
$ws.Range("E21").Value = 10.74942767388968
$ws.Range("G21").Value = 43.83816466327591
$ws.Range("H21").Value = 16.41271451984699
$ws.Range("M21").Value = 17.2812963929593
$ws.Range("N21").Value = 17.51918842261358
$ws.Range("B22").Value = 18.64572479733874
$ws.Range("C22").Value = 12.46437067893532
$ws.Range("D22").Value = 6.674901606269319
$ws.Range("E22").Value = 10.85264343565272
$ws.Range("G22").Value = 44.56417329700284
$ws.Range("H22").Value = 16.47846757808362
$ws.Range("M22").Value = 17.54529774097412
$ws.Range("N22").Value = 17.46829368602802
$ws.Range("B23").Value = 18.39041031837997
$ws.Range("C23").Value = 12.23914437259735
$ws.Range("D23").Value = 6.617100014778292
$ws.Range("E23").Value = 10.79740387530881
$ws.Range("G23").Value = 44.17689583961175
$ws.Range("H23").Value = 16.44303816786104
$ws.Range("M23").Value = 17.40451510759999
$ws.Range("N23").Value = 17.49530214891796
$ws.Range("B24").Value = 17.39655610711915
$ws.Range("C24").Value = 11.34538117075
$ws.Range("D24").Value = 6.39682902893631
$ws.Range("E24").Value = 10.59168227631287
$ws.Range("G24").Value = 42.70709462342155
$ws.Range("H24").Value = 16.31628908596119
$ws.Range("M24").Value = 16.86920962364304
$ws.Range("N24").Value = 17.60100059750459
$ws.Range("B25").Value = 16.27234246871258
$ws.Range("C25").Value = 10.29358546450806
$ws.Range("D25").Value = 6.15778237583713
$ws.Range("E25").Value = 10.37826245065687
$ws.Range("G25").Value = 41.12504834188655
$ws.Range("H25").Value = 16.19578548603395
$ws.Range("M25").Value = 16.2907044935307
$ws.Range("N25").Value = 17.72231264979244
